{"js": "// Remove author from metadata (docProps/core.xml dc:creator).\nconst properties = context.document.properties;\nproperties.author = \"\";\nawait context.sync();\n", "ps1": "# Remove author from metadata (docProps/core.xml dc:creator).\n$d = $word.ActiveDocument\n$d.Author = \"\"\n"}
